$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    # Row 2
    $ws.Cells.Item(2, 14).Value = 1.32  # N2
    $ws.Cells.Item(2, 16).Value = 1.32  # P2
    $ws.Cells.Item(2, 17).Value = 1.48  # Q2
    # Row 3
    $ws.Cells.Item(3, 7).Value = 1.37  # G3
    $ws.Cells.Item(3, 8).Value = 11  # H3
    $ws.Cells.Item(3, 20).Value = 1.91  # T3
    $ws.Cells.Item(3, 32).Value = 9.199999999999999  # AF3
    # Row 4
    $ws.Cells.Item(4, 7).Value = 5.7  # G4
    $ws.Cells.Item(4, 8).Value = 1.7  # H4
    $ws.Cells.Item(4, 9).Value = 1.72  # I4
    $ws.Cells.Item(4, 10).Value = 4.1  # J4
    $ws.Cells.Item(4, 15).Value = 1.23  # O4
    $ws.Cells.Item(4, 16).Value = 2.36  # P4
    $ws.Cells.Item(4, 17).Value = 1.7  # Q4
    $ws.Cells.Item(4, 18).Value = 1.55  # R4
    $ws.Cells.Item(4, 19).Value = 2.7  # S4
    $ws.Cells.Item(4, 27).Value = 17.5  # AA4
    $ws.Cells.Item(4, 29).Value = 9.6  # AC4
    $ws.Cells.Item(4, 32).Value = 46  # AF4
    $ws.Cells.Item(4, 34).Value = 17.5  # AH4
    $ws.Cells.Item(4, 36).Value = 150  # AJ4
    $ws.Cells.Item(4, 37).Value = 70  # AK4
    $ws.Cells.Item(4, 39).Value = 1000  # AM4
    $ws.Cells.Item(4, 41).Value = 7.8  # AO4
    # Row 6
    $ws.Cells.Item(6, 7).Value = 2.8  # G6
    $ws.Cells.Item(6, 8).Value = 2.64  # H6
    $ws.Cells.Item(6, 13).Value = 1.04  # M6
    $ws.Cells.Item(6, 14).Value = 5.2  # N6
    $ws.Cells.Item(6, 16).Value = 2.4  # P6
    $ws.Cells.Item(6, 17).Value = 1.69  # Q6
    # Row 7
    $ws.Cells.Item(7, 6).Value = 8.199999999999999  # F7
    $ws.Cells.Item(7, 7).Value = 9  # G7
    $ws.Cells.Item(7, 17).Value = 1.54  # Q7
    $ws.Cells.Item(7, 20).Value = 1.8  # T7
    $ws.Cells.Item(7, 24).Value = 28  # X7
    $ws.Cells.Item(7, 25).Value = 11.5  # Y7
    $ws.Cells.Item(7, 26).Value = 11  # Z7
    $ws.Cells.Item(7, 27).Value = 12.5  # AA7
    # Row 8
    $ws.Cells.Item(8, 6).Value = 2.24  # F8
    $ws.Cells.Item(8, 7).Value = 2.28  # G8
    $ws.Cells.Item(8, 8).Value = 3.4  # H8
    $ws.Cells.Item(8, 9).Value = 3.5  # I8
    $ws.Cells.Item(8, 13).Value = 1.05  # M8
    $ws.Cells.Item(8, 14).Value = 4.6  # N8
    $ws.Cells.Item(8, 15).Value = 1.25  # O8
    $ws.Cells.Item(8, 21).Value = 2.42  # U8
    $ws.Cells.Item(8, 24).Value = 18.5  # X8
    $ws.Cells.Item(8, 26).Value = 27  # Z8
    $ws.Cells.Item(8, 31).Value = 44  # AE8
    $ws.Cells.Item(8, 32).Value = 16  # AF8
    $ws.Cells.Item(8, 37).Value = 22  # AK8
    $ws.Cells.Item(8, 40).Value = 14  # AN8
    $ws.Cells.Item(8, 41).Value = 32  # AO8
    # Row 9
    $ws.Cells.Item(9, 6).Value = 1.7  # F9
    $ws.Cells.Item(9, 7).Value = 1.71  # G9
    $ws.Cells.Item(9, 8).Value = 5.1  # H9
    $ws.Cells.Item(9, 9).Value = 5.3  # I9
    $ws.Cells.Item(9, 10).Value = 4.4  # J9
    $ws.Cells.Item(9, 13).Value = 1.04  # M9
    $ws.Cells.Item(9, 17).Value = 1.6  # Q9
    $ws.Cells.Item(9, 20).Value = 1.64  # T9
    $ws.Cells.Item(9, 24).Value = 25  # X9
    $ws.Cells.Item(9, 26).Value = 130  # Z9
    $ws.Cells.Item(9, 29).Value = 10.5  # AC9
    $ws.Cells.Item(9, 32).Value = 13.5  # AF9
    $ws.Cells.Item(9, 34).Value = 17  # AH9
    $ws.Cells.Item(9, 41).Value = 160  # AO9
    # Row 10
    $ws.Cells.Item(10, 6).Value = 1.42  # F10
    $ws.Cells.Item(10, 8).Value = 8.800000000000001  # H10
    $ws.Cells.Item(10, 9).Value = 9.199999999999999  # I10
    $ws.Cells.Item(10, 16).Value = 2.68  # P10
    $ws.Cells.Item(10, 17).Value = 1.55  # Q10
    $ws.Cells.Item(10, 18).Value = 1.69  # R10
    $ws.Cells.Item(10, 20).Value = 1.81  # T10
    $ws.Cells.Item(10, 21).Value = 2.14  # U10
    $ws.Cells.Item(10, 24).Value = 28  # X10
    $ws.Cells.Item(10, 31).Value = 140  # AE10
    $ws.Cells.Item(10, 35).Value = 1000  # AI10
    $ws.Cells.Item(10, 36).Value = 12.5  # AJ10
    $ws.Cells.Item(10, 40).Value = 4.9  # AN10
    $ws.Cells.Item(10, 41).Value = 140  # AO10
    # Row 11
    $ws.Cells.Item(11, 11).Value = 6.8  # K11
    $ws.Cells.Item(11, 28).Value = 16  # AB11
    $ws.Cells.Item(11, 33).Value = 11  # AG11
    $ws.Cells.Item(11, 38).Value = 24  # AL11
    # Row 12
    $ws.Cells.Item(12, 7).Value = 2.42  # G12
    $ws.Cells.Item(12, 8).Value = 3.25  # H12
    $ws.Cells.Item(12, 9).Value = 3.35  # I12
    $ws.Cells.Item(12, 19).Value = 3.3  # S12
    # Row 13
    $ws.Cells.Item(13, 6).Value = 2.02  # F13
    $ws.Cells.Item(13, 9).Value = 4.4  # I13
    # Row 14
    $ws.Cells.Item(14, 6).Value = 1.8  # F14
    $ws.Cells.Item(14, 7).Value = 1.95  # G14
    $ws.Cells.Item(14, 8).Value = 4.3  # H14
    $ws.Cells.Item(14, 9).Value = 7  # I14
    $ws.Cells.Item(14, 11).Value = 4.1  # K14
    $ws.Cells.Item(14, 16).Value = 1.65  # P14
